# Apply updated "dSF" (column F) values from a data repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    4  = -2
    5  = -1
    6  = -8
    8  = -7
    9  = 1
    10 = 5
    12 = -2
    13 = 1
    15 = 6
    17 = -4
    19 = 5
    22 = -5
    23 = -3
    24 = 0
    26 = 0
    29 = -3
    30 = 0
    31 = 5
    32 = 1
    34 = -1
    36 = 0
    40 = -4
    43 = -3
    47 = 2
    48 = 0
    49 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
